# Extinction1.xlsx edit — "Reduced to 3 stim, added gray circles"
#
# 1) Rename the two existing stimulus labels (Extinction/... -> Stimuli/...)
# 2) Append a third stimulus block (rows 22-31) re-using the CS+4 label
# 3) Recolor the fonts to a consistent black and add a light-gray "gray
#    circle" fill banding to the CS+3 block and the new CS+4 block
# 4) Leave the active selection on E11

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$GRAY  = 15921906   # RGB(242,242,242) -> Excel "White, Background 1, Darker 5%"
$BLACK = 0          # RGB(0,0,0)

# ---- 1) Update the text of the first two stimulus blocks -----------------
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 1).Value = "Stimuli/CS+3.BMP"
}
for ($r = 12; $r -le 21; $r++) {
    $ws.Cells.Item($r, 1).Value = "Stimuli/CS+1.BMP"
}

# ---- 2) Styling of the existing blocks ------------------------------------
# Rows 12-21: unify the font color with column A (plain, no fill)
$ws.Range("A12:B21").Font.Color = $BLACK

# Rows 2-11 ("gray circle" banding): black font + light gray fill
$ws.Range("A2:B11").Font.Color = $BLACK
$ws.Range("A2:B11").Interior.Color = $GRAY

# ---- 3) Append the third stimulus block (rows 22-31) ----------------------
for ($r = 22; $r -le 31; $r++) {
    $ws.Cells.Item($r, 1).Value = "Stimuli/CS+4.BMP"
    $ws.Cells.Item($r, 2).Value = 0.4
}

# Copy the formatting (black font + gray fill) from the first banded block
# onto the new block so the same style objects are reused.
$ws.Range("A2:B11").Copy() | Out-Null
$ws.Range("A22:B31").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# ---- 4) Selection ----------------------------------------------------------
$ws.Range("E11").Select()
